$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08378199999999998
$ws.Range("N2").Value = 0.251346
$ws.Range("O2").Value = 0.007571394704126512
$ws.Range("P2").Value = 0.007571394704126512
$ws.Range("Q2").Value = 0.038238439928
$ws.Range("R2").Value = 0.344145959352
$ws.Range("S2").Value = 0.0001449892380990954
$ws.Range("T2").Value = 0.0001449892380990954

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("O3").Value = 0.006644889460697858
$ws.Range("P3").Value = 0.006644889460697857
$ws.Range("Q3").Value = 0.03355923398533334
$ws.Range("R3").Value = 0.302033105868
$ws.Range("S3").Value = 0.0001272470261831952
$ws.Range("T3").Value = 0.0001272470261831952

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.9857837158351757
$ws.Range("P4").Value = 0.9857837158351755
$ws.Range("Q4").Value = 4.97858490714
$ws.Range("R4").Value = 44.80726416426
$ws.Range("S4").Value = 0.01887737140576486
$ws.Range("T4").Value = 0.01887737140576486

# Row 5
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08378199999999998
$ws.Range("N5").Value = 0.251346
$ws.Range("O5").Value = 0.007571394704126512
$ws.Range("P5").Value = 0.007571394704126512
$ws.Range("Q5").Value = 1.654375486071333
$ws.Range("R5").Value = 14.889379374642
$ws.Range("S5").Value = 0.006272919128158824
$ws.Range("T5").Value = 0.006272919128158824

# Row 6
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("M6").Value = 0.07352966666666667
$ws.Range("O6").Value = 0.006644889460697858
$ws.Range("P6").Value = 0.006644889460697857
$ws.Range("S6").Value = 0.005505307255979515
$ws.Range("T6").Value = 0.005505307255979515

# Row 7
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.9857837158351757
$ws.Range("P7").Value = 0.9857837158351755
$ws.Range("S7").Value = 0.816724232316106
$ws.Range("T7").Value = 0.8167242323161059

# Row 8
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08378199999999998
$ws.Range("N8").Value = 0.251346
$ws.Range("O8").Value = 0.007571394704126512
$ws.Range("P8").Value = 0.007571394704126512
$ws.Range("Q8").Value = 0.304212358218
$ws.Range("R8").Value = 2.737911223962
$ws.Range("S8").Value = 0.001153486337868593
$ws.Range("T8").Value = 0.001153486337868593

# Row 9
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("M9").Value = 0.07352966666666667
$ws.Range("O9").Value = 0.006644889460697858
$ws.Range("P9").Value = 0.006644889460697857
$ws.Range("Q9").Value = 0.2669861461370001
$ws.Range("S9").Value = 0.001012335178535148
$ws.Range("T9").Value = 0.001012335178535147

# Row 10
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.9857837158351757
$ws.Range("P10").Value = 0.9857837158351755
$ws.Range("S10").Value = 0.1501821121133049
$ws.Range("T10").Value = 0.1501821121133049
